$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 15653.5
$ws.Range("I86").Value = 27675.75
$ws.Range("J86").Value = 3631.25
$ws.Range("K86").Value = 27675.75
$ws.Range("L86").Value = 3631.25
$ws.Range("M86").Value = -26552.75
$ws.Range("N86").Value = -5877.25

$ws.Range("H89").Value = 15653.5
$ws.Range("I89").Value = 27675.75
$ws.Range("J89").Value = 3631.25
$ws.Range("K89").Value = 138378.75
$ws.Range("L89").Value = 18156.25
$ws.Range("M89").Value = -132762.75
$ws.Range("N89").Value = -29388.25

$ws.Range("H137").Value = 5267889
$ws.Range("I137").Value = 6254780.5
$ws.Range("K137").Value = 18764341.5
$ws.Range("M137").Value = -18761791.5

$ws.Range("H141").Value = 917793.5600000001
$ws.Range("J141").Value = 4276170
$ws.Range("L141").Value = 12828510
$ws.Range("N141").Value = -12838870

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2420.625
$ws.Range("I2").Value = 1383
$ws.Range("K2").Value = 1383
$ws.Range("M2").Value = -1270

$ws.Range("H61").Value = 2272.9
$ws.Range("I61").Value = 1622.375
$ws.Range("J61").Value = 4875
$ws.Range("K61").Value = 1622.375
$ws.Range("L61").Value = 4875
$ws.Range("M61").Value = -1410.375
$ws.Range("N61").Value = -5299

$ws.Range("H74").Value = 842.125
$ws.Range("I74").Value = 817.5714
$ws.Range("K74").Value = 817.5714
$ws.Range("M74").Value = 56.42859999999996

$ws.Range("H77").Value = 842.125
$ws.Range("I77").Value = 817.5714
$ws.Range("K77").Value = 4087.857
$ws.Range("M77").Value = 280.143

$ws.Range("H102").Value = 2176.4348
$ws.Range("I102").Value = 2145.6191
$ws.Range("K102").Value = 2145.6191
$ws.Range("M102").Value = -523.6190999999999

$ws.Range("H116").Value = 2420.625
$ws.Range("I116").Value = 1383
$ws.Range("K116").Value = 1383
$ws.Range("M116").Value = 911

$ws.Range("H122").Value = 2814.2917
$ws.Range("I122").Value = 1863.5
$ws.Range("K122").Value = 5590.5
$ws.Range("M122").Value = -3140.5

$ws.Range("H132").Value = 23813794
$ws.Range("I132").Value = 33337320
$ws.Range("K132").Value = 100011960
$ws.Range("M132").Value = -100009430

$ws.Range("H136").Value = 2272.9
$ws.Range("I136").Value = 1622.375
$ws.Range("J136").Value = 4875
$ws.Range("K136").Value = 4867.125
$ws.Range("L136").Value = 14625
$ws.Range("M136").Value = -2317.125
$ws.Range("N136").Value = -19725

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2420.625
$ws.Range("I3").Value = 1383
$ws.Range("K3").Value = 1383
$ws.Range("M3").Value = -1269

$ws.Range("H107").Value = 2361.3635
$ws.Range("I107").Value = 1861.6
$ws.Range("K107").Value = 1861.6
$ws.Range("M107").Value = 58.40000000000009

$ws.Range("H134").Value = 6027
$ws.Range("I134").Value = 5945.909
$ws.Range("J134").Value = 6250
$ws.Range("K134").Value = 17837.727
$ws.Range("L134").Value = 18750
$ws.Range("M134").Value = -15302.727
$ws.Range("N134").Value = -23820

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3336674
$ws.Range("I31").Value = 4350509.5
$ws.Range("K31").Value = 4350509.5
$ws.Range("M31").Value = -4350214.5

$ws.Range("H34").Value = 3336674
$ws.Range("I34").Value = 4350509.5
$ws.Range("K34").Value = 4350509.5
$ws.Range("M34").Value = -4350307.5

$ws.Range("H58").Value = 31253740
$ws.Range("I58").Value = 1385
$ws.Range("J58").Value = 50005150
$ws.Range("K58").Value = 1385
$ws.Range("L58").Value = 50005150
$ws.Range("M58").Value = -1182
$ws.Range("N58").Value = -50005556

$ws.Range("H132").Value = 3058.2222
$ws.Range("I132").Value = 2369.5833
$ws.Range("J132").Value = 4435.5
$ws.Range("K132").Value = 7108.749899999999
$ws.Range("L132").Value = 13306.5
$ws.Range("M132").Value = -4578.749899999999
$ws.Range("N132").Value = -18366.5

$ws.Range("H134").Value = 1810.2
$ws.Range("I134").Value = 975.5625
$ws.Range("K134").Value = 2926.6875
$ws.Range("M134").Value = -391.6875

$ws.Range("H136").Value = 31253740
$ws.Range("I136").Value = 1385
$ws.Range("J136").Value = 50005150
$ws.Range("K136").Value = 4155
$ws.Range("L136").Value = 150015450
$ws.Range("M136").Value = -1605
$ws.Range("N136").Value = -150020550

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1875
$ws.Range("I5").Value = 568.5714
$ws.Range("J5").Value = 3399.1667
$ws.Range("K5").Value = 1705.7142
$ws.Range("L5").Value = 10197.5001
$ws.Range("M5").Value = -1593.7142
$ws.Range("N5").Value = -10421.5001

$ws.Range("H135").Value = 1875
$ws.Range("I135").Value = 568.5714
$ws.Range("J135").Value = 3399.1667
$ws.Range("K135").Value = 5117.1426
$ws.Range("L135").Value = 30592.5003
$ws.Range("M135").Value = -2582.1426
$ws.Range("N135").Value = -35662.5003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H101").Value = 34750
$ws.Range("J101").Value = 34750
$ws.Range("L101").Value = 34750
$ws.Range("N101").Value = -41240

$ws.Range("H102").Value = 2276.8823
$ws.Range("I102").Value = 1821.4166
$ws.Range("J102").Value = 3370
$ws.Range("K102").Value = 1821.4166
$ws.Range("L102").Value = 3370
$ws.Range("M102").Value = -199.4166
$ws.Range("N102").Value = -6614

$ws.Range("H132").Value = 2785.6743
$ws.Range("I132").Value = 2130.5862
$ws.Range("K132").Value = 6391.758600000001
$ws.Range("M132").Value = -3861.758600000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2015.5
$ws.Range("I93").Value = 1325.1538
$ws.Range("J93").Value = 2831.3635
$ws.Range("K93").Value = 1325.1538
$ws.Range("L93").Value = 2831.3635
$ws.Range("M93").Value = -77.15380000000005
$ws.Range("N93").Value = -5327.363499999999

$ws.Range("H127").Value = 31693.834
$ws.Range("J127").Value = 31693.834
$ws.Range("L127").Value = 31693.834
$ws.Range("N127").Value = -41613.834

$ws.Range("H132").Value = 3626.4348
$ws.Range("I132").Value = 1829.2858
$ws.Range("J132").Value = 4412.6875
$ws.Range("K132").Value = 5487.857400000001
$ws.Range("L132").Value = 13238.0625
$ws.Range("M132").Value = -2957.857400000001
$ws.Range("N132").Value = -18298.0625

$ws.Range("H136").Value = 3336013.8
$ws.Range("I136").Value = 4547692
$ws.Range("K136").Value = 13643076
$ws.Range("M136").Value = -13640526

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 1900.3334
$ws.Range("I62").Value = 1267.6666
$ws.Range("J62").Value = 2533
$ws.Range("K62").Value = 1267.6666
$ws.Range("L62").Value = 2533
$ws.Range("M62").Value = -643.6666
$ws.Range("N62").Value = -3781

$ws.Range("H65").Value = 1900.3334
$ws.Range("I65").Value = 1267.6666
$ws.Range("J65").Value = 2533
$ws.Range("K65").Value = 6338.333000000001
$ws.Range("L65").Value = 12665
$ws.Range("M65").Value = -3218.333000000001
$ws.Range("N65").Value = -18905

$ws.Range("H122").Value = 304706.34
$ws.Range("I122").Value = 358590.03
$ws.Range("K122").Value = 1075770.09
$ws.Range("M122").Value = -1073320.09

$ws.Range("H132").Value = 291935.38
$ws.Range("I132").Value = 402929.4
$ws.Range("J132").Value = 14450.3
$ws.Range("K132").Value = 1208788.2
$ws.Range("L132").Value = 43350.89999999999
$ws.Range("M132").Value = -1206258.2
$ws.Range("N132").Value = -48410.89999999999

$ws.Range("H136").Value = 1642.9
$ws.Range("I136").Value = 848.4286
$ws.Range("J136").Value = 3496.6667
$ws.Range("K136").Value = 2545.2858
$ws.Range("L136").Value = 10490.0001
$ws.Range("M136").Value = 4.714200000000346
$ws.Range("N136").Value = -15590.0001
